$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 197, pushing existing rows 197:245 down to 198:246
$ws.Rows("197:197").Insert()

# Populate the newly inserted row 197 with the new data point
$ws.Range("A197").Value = 10
$ws.Range("B197").Value = "Vega Modelo de Temuco"
$ws.Range("C197").Value = "La Araucanía"
$ws.Range("D197").Value = 44841
$ws.Range("E197").Value = 9
$ws.Range("F197").Value = 100112005
$ws.Range("G197").Value = "Puerro"
$ws.Range("H197").Value = "Azul de Maquehue"
$ws.Range("I197").Value = "Primera"
$ws.Range("J197").Value = 30
$ws.Range("K197").Value = 16000
$ws.Range("L197").Value = 16000
$ws.Range("M197").Value = 16000
$ws.Range("N197").Value = "$/docena de paquetes"
$ws.Range("O197").Value = "Provincia de Cautín"
$ws.Range("P197").Value = 1333
$ws.Range("Q197").Value = 12
$ws.Range("R197").Value = "Hortaliza"
